$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new text entries as described by the commit "se agrega imagen de demo en computacion y avances"
$ws.Range("A5").Value = "Para la ejecución de la solución presentada en el artículo se utilizaron las bases def2-SVP y los métodos BP86 y B3LYP, al consultar la documentación se encuentran todos disponibles, sin embargo una de las bases no se encuentra tal como la muestra el artículo, por lo que mientras se soluciona, se opta por utilizar en ambas la def2-SVP, la otra base era double-z def2-SVP"
$ws.Range("A6").Value = "Como se muestra en la anterior imágen el uso de orca para este sistema es simple en escritura, el paso posterior luego de comprobado que funciona la instalación, fue comunicar con el archivo en colab, sin embargo se encuentra en fase de despliegue del programa de orca para este proposito."

# Adjust row heights to match new content
$ws.Rows.Item(5).RowHeight = 60.75
$ws.Rows.Item(6).RowHeight = 45.75

# Update the view: scroll so row 4 is at top, and select A6
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A6").Select()
